# Actualización automática 2025-07-24 14:45:08
#
# Applies the recorded cell-value changes across the three worksheets:
#   - "VENTAS POR GRUPO"     (sheet1)
#   - "VENTA MENSUAL"        (sheet2)
#   - "CUMPLIMIENTO MENSUAL" (sheet3)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO"
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("I8").Value = 132.3
$wsGrupo.Range("M8").Value = 651.72

$wsGrupo.Range("D10").Value = 1465.34
$wsGrupo.Range("L10").Value = 1372.56

$wsGrupo.Range("D24").Value = "1 de 22"
$wsGrupo.Range("I24").Value = "1 de 22"
$wsGrupo.Range("L24").Value = "1 de 22"
$wsGrupo.Range("M24").Value = "5 de 22"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL"
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F8").Value = 784.02
$wsMensual.Range("F10").Value = 2837.9
$wsMensual.Range("F24").Value = 27081.55

# ---------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL"
# ---------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 3 - 240X80 PORCELANATO
$wsCumpl.Range("D3").Value = 1465.34
$wsCumpl.Range("E3").Value = 7203.57
$wsCumpl.Range("F3").Value = 0.1690339385228362

# Row 8 - LAVABOS
$wsCumpl.Range("D8").Value = 132.3
$wsCumpl.Range("E8").Value = 492.7
$wsCumpl.Range("F8").Value = 0.21168

# Row 15 - PIEDRA SINTERIZADA
$wsCumpl.Range("D15").Value = 1372.56
$wsCumpl.Range("E15").Value = 1128.45
$wsCumpl.Range("F15").Value = 0.5488022838773134

# Row 16 - PORCELANATO
$wsCumpl.Range("D16").Value = 24111.35
$wsCumpl.Range("E16").Value = 14645.19
$wsCumpl.Range("F16").Value = 0.6221233887235548

# Row 19 - TOTAL
$wsCumpl.Range("D19").Value = 27081.55
$wsCumpl.Range("E19").Value = 31141.45386304604
$wsCumpl.Range("F19").Value = 0.4651348814585738
